$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Danish term")
$ws.Columns("B").Delete()

# Delete old "(system)API Search Form" and "(system)API Search Form Sorting" columns
# After the deletion above, these are now at P and Q
$ws.Columns("P:Q").Delete()

# Insert a new column before O (currently "(system)API Search Criteria Mapping")
$ws.Columns("O").Insert()
$ws.Range("O1").Value = "(system)API Property Mapping"
$ws.Range("O2").Value = " "

# Row height for header
$ws.Rows("1").RowHeight = 14.25

# Autofit columns G:J
$ws.Range("G1:J2").EntireColumn.AutoFit()

# Selection
$ws.Range("A1:XFD1").Select()
